# Trade #115 closed at 2026-02-16 21:43:30 - leadlag UP +0.000%
#
# Trade #91 (leadlag, opened 21:38:25) is closed out with a time-exit, and a
# brand-new trade #115 (leadlag, UP) is opened at 21:43:30. This ripples
# through the per-strategy "leadlag" sheet, the combined "All Trades" log,
# and the two rollup sheets ("Summary", "Comparison").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": refresh the OVERALL / leadlag strategy rollup numbers.
# These text cells look numeric (percentages) so a leading "'" keeps Excel
# from re-interpreting them as percentage numbers, matching the original
# plain-text storage.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("C2").Value = 91
$wsSummary.Range("D2").Value = "'71.4%"
$wsSummary.Range("E2").Value = "'+27.2625%"
$wsSummary.Range("F2").Value = "'+0.2996%"

$wsSummary.Range("C3").Value = 88
$wsSummary.Range("D3").Value = "'48.9%"
$wsSummary.Range("E3").Value = "'+14.5771%"
$wsSummary.Range("F3").Value = "'+0.1656%"

# ---------------------------------------------------------------------------
# Sheet "leadlag": close out trade #91 (row 68) and append new trade #115
# (row 90).
# ---------------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")

$wsLeadlag.Range("G68").Value = 68112.23940400001
$wsLeadlag.Range("H68").Value = "CLOSED"
$wsLeadlag.Range("I68").Value = 0.4299
$wsLeadlag.Range("J68").Value = 4.3
$wsLeadlag.Range("M68").Value = "time_exit_5min"
$wsLeadlag.Range("N68").Value = 5

$wsLeadlag.Range("A90").Value = 115
$wsLeadlag.Range("B90").Value = "'2026-02-16"
$wsLeadlag.Range("C90").Value = "'21:43:30"
$wsLeadlag.Range("D90").Value = "leadlag"
$wsLeadlag.Range("E90").Value = "UP"
$wsLeadlag.Range("F90").Value = 68437.495
$wsLeadlag.Range("H90").Value = "OPEN"
$wsLeadlag.Range("I90").Value = 0
$wsLeadlag.Range("J90").Value = 0
$wsLeadlag.Range("K90").Value = 0.6051
$wsLeadlag.Range("L90").Value = "Binance leading with 0.061% move"
$wsLeadlag.Range("N90").Value = 0

# ---------------------------------------------------------------------------
# Sheet "All Trades": append the now-closed trade #91 as a new combined row
# (row 92).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Range("A92").Value = 91
$wsAll.Range("B92").Value = "'2026-02-16"
$wsAll.Range("C92").Value = "'21:38:25"
$wsAll.Range("D92").Value = "leadlag"
$wsAll.Range("E92").Value = "DOWN"
$wsAll.Range("F92").Value = 68406.285
$wsAll.Range("G92").Value = 68112.23940400001
$wsAll.Range("H92").Value = "CLOSED"
$wsAll.Range("I92").Value = 0.4299
$wsAll.Range("J92").Value = 4.3
$wsAll.Range("K92").Value = 0.75
$wsAll.Range("L92").Value = "Coinbase leading with -0.088% move"
$wsAll.Range("M92").Value = "time_exit_5min"
$wsAll.Range("N92").Value = 5

# ---------------------------------------------------------------------------
# Sheet "Comparison": refresh the leadlag strategy comparison row. D2/G2 look
# like plain numbers and C2/E2 look like percentages, so all four need the
# text-preserving "'" prefix.
# ---------------------------------------------------------------------------
$wsComparison = $wb.Worksheets.Item("Comparison")
$wsComparison.Range("B2").Value = 88
$wsComparison.Range("C2").Value = "'48.9%"
$wsComparison.Range("D2").Value = "'2.99"
$wsComparison.Range("E2").Value = "'+0.5093%"
$wsComparison.Range("G2").Value = "'1.67"
